$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.680.64"
$ws.Range("E2").Value = "  +0.16%  "

$ws.Range("D3").Value = "1.585.35"
$ws.Range("E3").Value = "  -1.93%  "

$ws.Range("E4").Value = "  +1.32%  "

$ws.Range("D5").Value = "206.32"
$ws.Range("E5").Value = "  -1.40%  "

$ws.Range("E6").Value = "  -1.82%  "

$ws.Range("E7").Value = "  +1.35%  "

$ws.Range("D8").Value = "22.23"
$ws.Range("E8").Value = "  -3.47%  "

$ws.Range("E9").Value = "  -0.54%  "

$ws.Range("E10").Value = "  -2.50%  "

$ws.Range("D11").Value = "0.0865"
$ws.Range("E11").Value = "  -0.87%  "

$ws.Range("E12").Value = "  -2.00%  "

$ws.Range("D13").Value = "1.585.64"
$ws.Range("E13").Value = "  -3.11%  "

$ws.Range("E14").Value = "  -2.36%  "

$ws.Range("D15").Value = "0.532"
$ws.Range("E15").Value = "  -4.29%  "

$ws.Range("D16").Value = "27.652.12"
$ws.Range("E16").Value = "  +0.00%  "

$ws.Range("D17").Value = "63.29"
$ws.Range("E17").Value = "  -2.09%  "

$ws.Range("D18").Value = "219.52"
$ws.Range("E18").Value = "  -3.33%  "

$ws.Range("D19").Value = "0.0₃0693"
$ws.Range("E19").Value = "  -2.82%  "

$ws.Range("D20").Value = "7.31"
$ws.Range("E20").Value = "  -4.00%  "

$ws.Range("E21").Value = "  +1.42%  "

$ws.Range("D22").Value = "4.15"
$ws.Range("E22").Value = "  -3.51%  "

$ws.Range("D23").Value = "9.56"
$ws.Range("E23").Value = "  -4.61%  "

$ws.Range("E24").Value = "  -2.28%  "

$ws.Range("D25").Value = "155.15"
$ws.Range("E25").Value = "  +0.44%  "

$ws.Range("D26").Value = "6.84"
$ws.Range("E26").Value = "  -0.48%  "

$ws.Range("E27").Value = "  +1.36%  "

$ws.Range("D28").Value = "15.11"
$ws.Range("E28").Value = "  -1.90%  "

$ws.Range("E29").Value = "  -3.27%  "

$ws.Range("E30").Value = "  -1.48%  "

$ws.Range("D31").Value = "0.0467"
$ws.Range("E31").Value = "  -2.27%  "

$ws.Range("E32").Value = "  -3.22%  "

$ws.Range("D33").Value = "1.381.34"
$ws.Range("E33").Value = "  -0.53%  "

$ws.Range("D34").Value = "2.94"
$ws.Range("E34").Value = "  -4.14%  "

$ws.Range("E35").Value = "  -3.09%  "

$ws.Range("D36").Value = "0.980"
$ws.Range("E36").Value = "  -1.48%  "

$ws.Range("E38").Value = "  -2.62%  "

$ws.Range("D39").Value = "0.539"
$ws.Range("E39").Value = "  -2.87%  "

$ws.Range("E40").Value = "  -2.03%  "

$ws.Range("E41").Value = "  +1.39%  "

$ws.Range("E42").Value = "  -2.79%  "

$ws.Range("E43").Value = "  -2.54%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "63.50"
$ws.Range("E44").Value = "  -2.91%  "

$ws.Range("B45").Value = "MXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D45").Value = "2.17"
$ws.Range("E45").Value = "  +2.73%  "

$ws.Range("E46").Value = "  -2.45%  "

$ws.Range("D47").Value = "1.721.07"
$ws.Range("E47").Value = "  -1.96%  "

$ws.Range("E49").Value = "  +10.54%  "

$ws.Range("D50").Value = "0.0972"
$ws.Range("E50").Value = "  -3.41%  "

$ws.Range("E51").Value = "  -0.56%  "
